# Weekly update: insert a new "Sandia" (watermelon) price record as the new
# row 234 on the "Vega Modelo de Temuco" sheet, pushing all subsequent
# records down by one row (old row 234 becomes 235, ..., old row 342
# becomes the new row 343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 234; everything below shifts down one row
# (this also grows the used range / dimension from R342 to R343).
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record.
$ws.Cells.Item(234, 1).Value  = 10
$ws.Cells.Item(234, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(234, 3).Value  = "La Araucanía"
$ws.Cells.Item(234, 4).Value  = 44523
$ws.Cells.Item(234, 5).Value  = 9
$ws.Cells.Item(234, 6).Value  = 100112028
$ws.Cells.Item(234, 7).Value  = "Sandia"
$ws.Cells.Item(234, 8).Value  = "Sin especificar"
$ws.Cells.Item(234, 9).Value  = "Primera"
$ws.Cells.Item(234, 10).Value = 250
$ws.Cells.Item(234, 11).Value = 950
$ws.Cells.Item(234, 12).Value = 950
$ws.Cells.Item(234, 13).Value = 950
$ws.Cells.Item(234, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(234, 15).Value = "Perú"
$ws.Cells.Item(234, 16).Value = 950
$ws.Cells.Item(234, 17).Value = 1
$ws.Cells.Item(234, 18).Value = "Hortaliza"
